$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Remove the old rows 4-14 (and their header/body no longer needed)
#    Clearing the range shrinks the sheet's used range / dimension to
#    A1:H3, matching the target (only header + 2 data rows remain).
# ------------------------------------------------------------------
$ws.Range("A4:H14").Clear()

# ------------------------------------------------------------------
# 2. Drop every existing hyperlink (there were 13, one per old data
#    row) - this engine's Hyperlinks.Delete() clears them workbook/
#    sheet-wide regardless of the range it is invoked from, so a
#    single call is enough to reset to zero hyperlinks.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 3. Row 2: refresh timestamp + replace with the "airline matching
#    app" listing (previously row 9, now promoted to row 2 with an
#    updated scrape timestamp).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "2025-12-26 06:29:28"
$ws.Range("B2").Value = "【急募】航空会社とお客様のマッチングサービスのアプリ開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5461280"
$ws.Range("G2").Value = 85
$ws.Range("H2").Value = "◆開発 ◇アプリ"

# ------------------------------------------------------------------
# 4. Row 3: refresh timestamp + replace with the brand-new "bidding
#    DX cloud MVP" listing.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "2025-12-26 06:29:28"
$ws.Range("B3").Value = "【急募】入札DXクラウドMVP開発者募集!"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5461481"
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = "◆開発"

# ------------------------------------------------------------------
# 5. Re-create the two hyperlinks that should remain, on F2 and F3,
#    pointing at their respective URLs. Re-apply the named Hyperlink
#    style afterwards so the cell keeps using the workbook's original
#    "Hyperlink" cell-style record instead of a freshly synthesised
#    (but equivalent) one that Hyperlinks.Add() tends to allocate.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5461280")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5461481")
$ws.Range("F3").Style = "Hyperlink"

# ------------------------------------------------------------------
# 6. Column width tweaks (B 52->30, D 30->28, H 17->12). This engine
#    stores width as ColumnWidth + 0.8333333333333334, so back that
#    offset out of the target widths to land on the exact values.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 29.166666666666668
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
